# Auto-generated edit script applying numeric cell updates described by the
# commit diff to Sheets/Golem_Profits.xlsx (workbook sheets ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Only plain numeric cell values change; no formulas or
# formatting are touched.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 304.625
$ws.Range("I15").Value = 304.625
$ws.Range("K15").Value = 913.875
$ws.Range("M15").Value = -744.875
$ws.Range("H34").Value = 4050
$ws.Range("I34").Value = 3575.25
$ws.Range("K34").Value = 3575.25
$ws.Range("M34").Value = -3372.25
$ws.Range("H36").Value = 4050
$ws.Range("I36").Value = 3575.25
$ws.Range("K36").Value = 3575.25
$ws.Range("M36").Value = -2860.25
$ws.Range("H62").Value = 999.5
$ws.Range("J62").Value = 999
$ws.Range("L62").Value = 999
$ws.Range("N62").Value = -2247
$ws.Range("H65").Value = 999.5
$ws.Range("J65").Value = 999
$ws.Range("L65").Value = 4995
$ws.Range("N65").Value = -11235
$ws.Range("H74").Value = 5004
$ws.Range("J74").Value = 5004
$ws.Range("L74").Value = 5004
$ws.Range("N74").Value = -6876
$ws.Range("H77").Value = 5004
$ws.Range("J77").Value = 5004
$ws.Range("L77").Value = 25020
$ws.Range("N77").Value = -34380
$ws.Range("H88").Value = 1000
$ws.Range("J88").Value = 1000
$ws.Range("L88").Value = 1000
$ws.Range("N88").Value = -1812
$ws.Range("H91").Value = 1000
$ws.Range("J91").Value = 1000
$ws.Range("L91").Value = 1000
$ws.Range("N91").Value = -3808
$ws.Range("H92").Value = 969.8
$ws.Range("I92").Value = 1399.6666
$ws.Range("K92").Value = 1399.6666
$ws.Range("M92").Value = -151.6666
$ws.Range("H98").Value = 15167.75
$ws.Range("I98").Value = 6939.6
$ws.Range("J98").Value = 28881.334
$ws.Range("K98").Value = 6939.6
$ws.Range("L98").Value = 28881.334
$ws.Range("M98").Value = -5441.6
$ws.Range("N98").Value = -31877.334
$ws.Range("H122").Value = 15167.75
$ws.Range("I122").Value = 6939.6
$ws.Range("J122").Value = 28881.334
$ws.Range("K122").Value = 20818.8
$ws.Range("L122").Value = 86644.00199999999
$ws.Range("M122").Value = -18368.8
$ws.Range("N122").Value = -91544.00199999999
$ws.Range("H137").Value = 937.2222
$ws.Range("I137").Value = 947.8570999999999
$ws.Range("K137").Value = 2843.5713
$ws.Range("M137").Value = -293.5712999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1664.4615
$ws.Range("I61").Value = 1557.3636
$ws.Range("K61").Value = 1557.3636
$ws.Range("M61").Value = -1345.3636
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050
$ws.Range("H132").Value = 1900
$ws.Range("I132").Value = 1900
$ws.Range("K132").Value = 5700
$ws.Range("M132").Value = -3170
$ws.Range("H136").Value = 1664.4615
$ws.Range("I136").Value = 1557.3636
$ws.Range("K136").Value = 4672.0908
$ws.Range("M136").Value = -2122.0908

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 20000
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20620
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4997.5
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 4995
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 4995
$ws.Range("M4").Value = -4888
$ws.Range("N4").Value = -5219
$ws.Range("H19").Value = 187.3077
$ws.Range("I19").Value = 187.3077
$ws.Range("K19").Value = 187.3077
$ws.Range("M19").Value = -17.30770000000001
$ws.Range("H24").Value = 187.3077
$ws.Range("I24").Value = 187.3077
$ws.Range("K24").Value = 187.3077
$ws.Range("M24").Value = -17.30770000000001
$ws.Range("H33").Value = 6234.8335
$ws.Range("J33").Value = 21999.666
$ws.Range("L33").Value = 21999.666
$ws.Range("N33").Value = -22757.666
$ws.Range("H35").Value = 5499.5
$ws.Range("J35").Value = 9999
$ws.Range("L35").Value = 9999
$ws.Range("N35").Value = -10587
$ws.Range("H36").Value = 4999
$ws.Range("I36").Value = 5500
$ws.Range("J36").Value = 3496
$ws.Range("K36").Value = 5500
$ws.Range("L36").Value = 3496
$ws.Range("M36").Value = -5112
$ws.Range("N36").Value = -4272
$ws.Range("H40").Value = 4999
$ws.Range("I40").Value = 5500
$ws.Range("J40").Value = 3496
$ws.Range("K40").Value = 5500
$ws.Range("L40").Value = 3496
$ws.Range("M40").Value = -5340
$ws.Range("N40").Value = -3816
$ws.Range("H42").Value = 26666
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 26666
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 26666
$ws.Range("M42").Value = $null
$ws.Range("N42").Value = -27852
$ws.Range("H43").Value = 45928.332
$ws.Range("J43").Value = 45928.332
$ws.Range("L43").Value = 45928.332
$ws.Range("N43").Value = -46296.332
$ws.Range("H44").Value = 29998.6
$ws.Range("I44").Value = 29997.666
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 29997.666
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = -29555.666
$ws.Range("N44").Value = -30884
$ws.Range("H58").Value = 1532
$ws.Range("I58").Value = 1456.8572
$ws.Range("K58").Value = 1456.8572
$ws.Range("M58").Value = -1253.8572
$ws.Range("H96").Value = 30517.572
$ws.Range("J96").Value = 34603.832
$ws.Range("L96").Value = 34603.832
$ws.Range("N96").Value = -40095.832
$ws.Range("H101").Value = 45928.332
$ws.Range("J101").Value = 45928.332
$ws.Range("L101").Value = 45928.332
$ws.Range("N101").Value = -52418.332
$ws.Range("H103").Value = 20000
$ws.Range("I103").Value = 20000
$ws.Range("K103").Value = 20000
$ws.Range("M103").Value = -18828
$ws.Range("H136").Value = 1532
$ws.Range("I136").Value = 1456.8572
$ws.Range("K136").Value = 4370.571599999999
$ws.Range("M136").Value = -1820.571599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 20.6
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 29.2
$ws.Range("K2").Value = 72
$ws.Range("L2").Value = 175.2
$ws.Range("M2").Value = 41
$ws.Range("N2").Value = -401.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 1250
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 1500
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = -860
$ws.Range("N12").Value = -1780
$ws.Range("H31").Value = 1569
$ws.Range("I31").Value = 1569
$ws.Range("K31").Value = 1569
$ws.Range("M31").Value = -1277
$ws.Range("H37").Value = 1569
$ws.Range("I37").Value = 1569
$ws.Range("K37").Value = 1569
$ws.Range("M37").Value = -1292
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H126").Value = 12497.875
$ws.Range("I126").Value = 9996.6
$ws.Range("K126").Value = 29989.8
$ws.Range("M126").Value = -27519.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 341.66666
$ws.Range("I9").Value = 325
$ws.Range("J9").Value = 375
$ws.Range("K9").Value = 325
$ws.Range("L9").Value = 375
$ws.Range("M9").Value = -101
$ws.Range("N9").Value = -823
$ws.Range("H35").Value = 4199
$ws.Range("I35").Value = 4192.5
$ws.Range("J35").Value = 4208.75
$ws.Range("K35").Value = 4192.5
$ws.Range("L35").Value = 4208.75
$ws.Range("M35").Value = -3856.5
$ws.Range("N35").Value = -4880.75
$ws.Range("H40").Value = 32636.363
$ws.Range("I40").Value = 21124.5
$ws.Range("K40").Value = 21124.5
$ws.Range("M40").Value = -20988.5
$ws.Range("H50").Value = 30084
$ws.Range("J50").Value = 30084
$ws.Range("L50").Value = 30084
$ws.Range("N50").Value = -31358
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11496.728
$ws.Range("I136").Value = 9294.4
$ws.Range("J136").Value = 13332
$ws.Range("K136").Value = 27883.2
$ws.Range("L136").Value = 39996
$ws.Range("M136").Value = -25333.2
$ws.Range("N136").Value = -45096

